$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row (row 11): Right marks per question 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row (row 12): Right total 57 -> 95
$ws.Range("B12").Value = 95

# Update the Corr/total marks text 56/84 -> 95/140
$ws.Range("E12").Value = "95/140"
